$d = $word.ActiveDocument

# Helper: find the paragraph whose text contains $needle (first match).
function Find-ParagraphContaining($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text -like ("*" + $needle + "*")) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "The generation of code documentation of the backend part requires an" +
#    " " + "execution of" -> touch/normalize into a single run (the resulting
#    text is unchanged, only the run split collapses).
# ---------------------------------------------------------------------------
$pBackendRequires = Find-ParagraphContaining $d "requires an execution of"
$pBackendRequires.Range.Find.Execute(
    "The generation of code documentation of the backend part requires an execution of",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "The generation of code documentation of the backend part requires an execution of",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Fix the copy/paste mistake: the *second* "jsdoc -c jsdocConfigFrontend.json"
#    quoted command (the one right after the backend-specific instructions)
#    should read "...jsdocConfigBackend.json".
# ---------------------------------------------------------------------------
$pBackendCmd = Find-ParagraphContaining $d "jsdocConfigFrontend.json"
$prevText = ""
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*jsdocConfigFrontend.json*") {
        $prev = $d.Paragraphs.Item($i - 1)
        if ($prev.Range.Text -like "*requires an execution of*") {
            $pBackendCmd = $p
        }
    }
}
$pBackendCmd.Range.Find.Execute("Frontend", $true, $false, $false, $false, $false, $true, 1, $false, "Backend", 2) | Out-Null

# Move the "_GoBack" bookmark so it sits right after the "Backend" word that
# was just fixed (Word only ever keeps one "_GoBack" bookmark; adding a new
# one with the same name re-homes it, removing it from its previous spot at
# the very end of the document).
$pBackendCmd = Find-ParagraphContaining $d "jsdocConfigBackend.json"
$cmdText = $pBackendCmd.Range.Text
$idx = $cmdText.IndexOf("Backend")
$bmStart = $pBackendCmd.Range.Start + $idx + "Backend".Length
$bm = $d.Range($bmStart, $bmStart)
$d.Bookmarks.Add("_GoBack", $bm) | Out-Null

# ---------------------------------------------------------------------------
# 3) Grammar fix: "After the execution aforementioned commands..." should
#    read "After the execution of aforementioned commands...".
# ---------------------------------------------------------------------------
$pAfterExecution = Find-ParagraphContaining $d "After the execution"
$pAfterExecution.Range.Find.Execute(
    "After the execution aforementioned ",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "After the execution of aforementioned ",
    2) | Out-Null
